$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.02953839302063
$ws.Range("B1").Value = 2.338687181472778
$ws.Range("C1").Value = 4.865468978881836
$ws.Range("D1").Value = 2.460928440093994
$ws.Range("E1").Value = 1.338848114013672
